$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column after the existing header row (H1), copying the
# formatting (bold + border + centered alignment) already used by the other
# header cells such as G1 ("sum").
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the numeric "Save" value for the data row.
$ws.Range("H2").Value = 1
